# Applies proofing-error markup (w:proofErr) splits to three list
# paragraphs and removes the "Gestão de promotores/funcionários ..."
# bullet entirely (leaving an empty paragraph behind), matching the
# "Base de Dados a funcionar" commit.

$d = $word.ActiveDocument

function Set-ParagraphInnerXml($paragraph, [string]$innerXml) {
    # Replace only the paragraph's *content* (everything up to, but not
    # including, its trailing paragraph mark) so the paragraph's own
    # <w:pPr> (style/numbering) survives untouched.
    $full = $paragraph.Range
    $target = $d.Range($full.Start, $full.End - 1)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xml)
}

# --- Paragraph: "Pagina login - em que só o adm pode adicionar o novo funcionário " ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Pagina login - em que")) {
        $inner = '<w:proofErr w:type="gramStart"/>' +
                 '<w:r><w:t>Pagina</w:t></w:r>' +
                 '<w:proofErr w:type="gramEnd"/>' +
                 '<w:r><w:t xml:space="preserve"> login - em que só o </w:t></w:r>' +
                 '<w:proofErr w:type="spellStart"/>' +
                 '<w:r><w:t>adm</w:t></w:r>' +
                 '<w:proofErr w:type="spellEnd"/>' +
                 '<w:r><w:t xml:space="preserve"> pode adicionar o novo funcionário </w:t></w:r>'
        Set-ParagraphInnerXml $p $inner
        break
    }
}

# --- Paragraph: "Gestão de promotores/funcionários ..." -> removed entirely ---
# Deleting the *whole* paragraph Range (including its paragraph mark)
# merges it away, leaving the following (previously empty) paragraph's
# mark in place with no leftover <w:pPr>.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Gestão de promotores/funcionários")) {
        $p.Range.Delete()
        break
    }
}

# --- Paragraph: "Admin adiciona novo funcionário (programa devolve palavra passe)" ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Admin adiciona novo funcionário")) {
        $inner = '<w:proofErr w:type="spellStart"/>' +
                 '<w:r><w:t>Admin</w:t></w:r>' +
                 '<w:proofErr w:type="spellEnd"/>' +
                 '<w:r><w:t xml:space="preserve"> adiciona novo funcionário (programa devolve </w:t></w:r>' +
                 '<w:proofErr w:type="gramStart"/>' +
                 '<w:r><w:t>palavra passe</w:t></w:r>' +
                 '<w:proofErr w:type="gramEnd"/>' +
                 '<w:r><w:t>)</w:t></w:r>'
        Set-ParagraphInnerXml $p $inner
        break
    }
}

# --- Paragraph: "Funcionário entra com essa palavra passe, o programa pede ..." ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Funcionário entra com essa")) {
        $inner = '<w:r><w:t xml:space="preserve">Funcionário entra com essa </w:t></w:r>' +
                 '<w:proofErr w:type="gramStart"/>' +
                 '<w:r><w:t>palavra passe</w:t></w:r>' +
                 '<w:proofErr w:type="gramEnd"/>' +
                 '<w:r><w:t>, o programa pede para alterar a palavra passe</w:t></w:r>'
        Set-ParagraphInnerXml $p $inner
        break
    }
}

Write-Output "done"
